# Auto-generated COM-interop script that reproduces the "Add git concept
# about submodule and staging area" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('A22').Value = 'Git concept'
$ws.Range('B22').Value = 'Relationship between "working directory" and "staging area"'
$ws.Range('A23').Value = 'Git concept'
$ws.Range('B23').Value = 'Git folder under git folder: submodule'

$ws.Range('C22').Value = 'STAGE 本質上 is the product of comparasion between HEAD and WORKDIR. So items in STAGE can be regarded as " actions" basic on HEAD to become WORK. Below command or action is related:
> git add {…}      //move red to green entry in git status that add stuff to head
> git rm {…}       //move red to green entry in git status that remove stuffs from head
> git rm --cached {...}       //move green to red to unstage
> git checkout -- {…}      //move red to nth in git status
> user edit the WORK      //generate red entry in git status'
$stage0 = $ws.Range('C22').Characters(1, 206)
$stage0.Font.Size = 10
$stage0.Font.Color = 0
$stage1 = $ws.Range('C22').Characters(207, 3)
$stage1.Font.Size = 10
$stage1.Font.Color = 255
$stage2 = $ws.Range('C22').Characters(210, 4)
$stage2.Font.Size = 10
$stage2.Font.Color = 0
$stage3 = $ws.Range('C22').Characters(214, 5)
$stage3.Font.Size = 10
$stage3.Font.Color = 5287936
$stage4 = $ws.Range('C22').Characters(219, 70)
$stage4.Font.Size = 10
$stage4.Font.Color = 0
$stage5 = $ws.Range('C22').Characters(289, 3)
$stage5.Font.Size = 10
$stage5.Font.Color = 255
$stage6 = $ws.Range('C22').Characters(292, 4)
$stage6.Font.Size = 10
$stage6.Font.Color = 0
$stage7 = $ws.Range('C22').Characters(296, 5)
$stage7.Font.Size = 10
$stage7.Font.Color = 5287936
$stage8 = $ws.Range('C22').Characters(301, 87)
$stage8.Font.Size = 10
$stage8.Font.Color = 0
$stage9 = $ws.Range('C22').Characters(388, 5)
$stage9.Font.Size = 10
$stage9.Font.Color = 5287936
$stage10 = $ws.Range('C22').Characters(393, 4)
$stage10.Font.Size = 10
$stage10.Font.Color = 0
$stage11 = $ws.Range('C22').Characters(397, 3)
$stage11.Font.Size = 10
$stage11.Font.Color = 255
$stage12 = $ws.Range('C22').Characters(400, 46)
$stage12.Font.Size = 10
$stage12.Font.Color = 0
$stage13 = $ws.Range('C22').Characters(446, 3)
$stage13.Font.Size = 10
$stage13.Font.Color = 255
$stage14 = $ws.Range('C22').Characters(449, 59)
$stage14.Font.Size = 10
$stage14.Font.Color = 0
$stage15 = $ws.Range('C22').Characters(508, 3)
$stage15.Font.Size = 10
$stage15.Font.Color = 255
$stage16 = $ws.Range('C22').Characters(511, 20)
$stage16.Font.Size = 10
$stage16.Font.Color = 0

$ws.Range('C23').Value = '* When our new git folder contains another git folder, then the subfolder is name "submodule" in git machanism.
* When do adding in parent git folder, submodule will present as {path}/{submodule_name} as a standalone item, then no details inside submodule will be revealed in the parent folder git console.
* When submodule has some uncommited change (staged or not), the parent folder git status will note the changes in red item as below:
    > modified:   {path}/SUBMODULE1 (modified content, untracked content)
    > modified:   {path}/SUBMODULE2 (modified content)
And user cannot stage above "changes" unless going inside submodule and clean the directory.
(p.s. the error message when try to add the submodule in parent like: Git: fatal: Pathspec is in submodule )'
$sub0 = $ws.Range('C23').Characters(1, 733)
$sub0.Font.Size = 10
$sub0.Font.Color = 0
$sub1 = $ws.Range('C23').Characters(734, 36)
$sub1.Font.Size = 10
$sub1.Font.Italic = $true
$sub1.Font.Color = 8421504
$sub2 = $ws.Range('C23').Characters(770, 2)
$sub2.Font.Size = 10
$sub2.Font.Color = 0

$ws.Rows.Item(22).RowHeight = 90
$ws.Rows.Item(23).RowHeight = 141

$ws.Range('C21').Select()

Write-Output "edit applied"
